# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# Updates the worker account-statement rows (B16:J27) on sheet "Hoja1":
# each worker's document id / name / period / overdue amount ("Valor Mora")
# / base salary ("Salario Basico") is refreshed with the new database values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: C = N Doc Trabajador, D = Nombre Trabajador, E = Periodo Mora,
#          F = Valor Mora, G = Salario Basico

$data = @(
    @{ Row = 16; Doc = "1051448153"; Nombre = "AMIR PAJARO PAJARO";             Periodo = "2303"; Mora = 46400; Salario = 1160000 },
    @{ Row = 17; Doc = "1051448153"; Nombre = "AMIR PAJARO PAJARO";             Periodo = "2210"; Mora = 40000; Salario = 1160000 },
    @{ Row = 18; Doc = "1001898501"; Nombre = "DIDIER ANDRES MIRANDA SANCHEZ";  Periodo = "2303"; Mora = 46400; Salario = 1160000 },
    @{ Row = 19; Doc = "1001898501"; Nombre = "DIDIER ANDRES MIRANDA SANCHEZ";  Periodo = "2302"; Mora = 40000; Salario = 1160000 },
    @{ Row = 20; Doc = "1001898501"; Nombre = "DIDIER ANDRES MIRANDA SANCHEZ";  Periodo = "2210"; Mora = 40000; Salario = 1160000 },
    @{ Row = 21; Doc = "1007856972"; Nombre = "HYLEANA MARGARITA BARRIOS PUERTA"; Periodo = "2304"; Mora = 32707; Salario = 908526 },
    @{ Row = 22; Doc = "1007856972"; Nombre = "HYLEANA MARGARITA BARRIOS PUERTA"; Periodo = "2303"; Mora = 46400; Salario = 908526 },
    @{ Row = 23; Doc = "1007856972"; Nombre = "HYLEANA MARGARITA BARRIOS PUERTA"; Periodo = "2302"; Mora = 46400; Salario = 908526 },
    @{ Row = 24; Doc = "1007856972"; Nombre = "HYLEANA MARGARITA BARRIOS PUERTA"; Periodo = "2210"; Mora = 40000; Salario = 908526 },
    @{ Row = 25; Doc = "1099962566"; Nombre = "MILTON JESUS CONDE LOZANO";      Periodo = "2303"; Mora = 46400; Salario = 1160000 },
    @{ Row = 26; Doc = "1099962566"; Nombre = "MILTON JESUS CONDE LOZANO";      Periodo = "2302"; Mora = 46400; Salario = 1160000 },
    @{ Row = 27; Doc = "1099962566"; Nombre = "MILTON JESUS CONDE LOZANO";      Periodo = "2210"; Mora = 40000; Salario = 1160000 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Mora
    $ws.Cells.Item($r, 7).Value = $item.Salario
}

$wb.Save()
